$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cell values (minor corrections to existing figures) ---
$ws.Range("Y57").Value = 8456794.624
$ws.Range("AC57").Value = 8642083.84
$ws.Range("AG57").Value = 9867097.088
$ws.Range("AO57").Value = 12515202.048
$ws.Range("Y58").Value = -3953172.736
$ws.Range("AC58").Value = -3993735.68
$ws.Range("AK58").Value = -5648540.16
$ws.Range("U59").Value = 2809900.8
$ws.Range("AC59").Value = 4648346.624
$ws.Range("AK59").Value = 5826446.336
$ws.Range("BA59").Value = 8378442.24
$ws.Range("BE59").Value = 10133890.048
$ws.Range("AG60").Value = -1795960.704
$ws.Range("AK60").Value = -1887737.984
$ws.Range("AO60").Value = -2201873.152
$ws.Range("BE60").Value = -3056449.28
$ws.Range("AG61").Value = 3473515.008
$ws.Range("AK61").Value = 3938708.736
$ws.Range("AO61").Value = 4300455.424
$ws.Range("BE61").Value = 7077436.928
$ws.Range("U62").Value = -534396.96
$ws.Range("AG62").Value = -1057344.896
$ws.Range("BE62").Value = -1974976.256
$ws.Range("AC63").Value = -518061.984
$ws.Range("AG63").Value = -31423.968
$ws.Range("AW63").Value = -481449.056
$ws.Range("BA63").Value = -369941.024
$ws.Range("BE63").Value = -395064.032
$ws.Range("Y65").Value = 271742.976
$ws.Range("AC65").Value = 31853.976
$ws.Range("AG65").Value = 112683.992
$ws.Range("AS65").Value = -47864.992
$ws.Range("U66").Value = -69702.992
$ws.Range("AC66").Value = -202404.032
$ws.Range("AS67").Value = -1280946.944
$ws.Range("Y68").Value = -220592.976
$ws.Range("AC68").Value = -318374.944
$ws.Range("AO68").Value = -199520.032
$ws.Range("BA68").Value = -90981
$ws.Range("BE68").Value = -240226.016
$ws.Range("U69").Value = 46517.032
$ws.Range("Y69").Value = -57787.984
$ws.Range("AW69").Value = 291058.016
$ws.Range("BE69").Value = 146062.032
$ws.Range("U70").Value = -250803.984
$ws.Range("AK70").Value = -348231.008
$ws.Range("AO70").Value = -468635.072
$ws.Range("AS70").Value = 108770.976
$ws.Range("BA70").Value = -239116.96
$ws.Range("AC71").Value = -60903.992
$ws.Range("AC73").Value = -47707.016
$ws.Range("U74").Value = 466145.984
$ws.Range("Y74").Value = 799721.9840000001
$ws.Range("AO74").Value = 1510798.08
$ws.Range("AS74").Value = 3749697.28
$ws.Range("BE74").Value = 4757251.584
$ws.Range("Y75").Value = -506704.992
$ws.Range("BE75").Value = -411236.064
$ws.Range("U76").Value = 44191
$ws.Range("AK76").Value = -142743.968
$ws.Range("AC80").Value = 697188.032
$ws.Range("AO80").Value = 964534.976

# --- Clear cells that should become blank (previously zero placeholders) ---
$ws.Range("AW57:BH57").ClearContents()
$ws.Range("AW58:BH58").ClearContents()
$ws.Range("S64:AW64").ClearContents()
$ws.Range("AW71:BH71").ClearContents()
$ws.Range("AW72:BH72").ClearContents()
$ws.Range("AW73:BH73").ClearContents()
$ws.Range("AW77:BH77").ClearContents()
$ws.Range("AW78:BH78").ClearContents()
$ws.Range("AP79:AW79").ClearContents()
